# Generate Report for Handback
# Adds a new handback record (01504b35-5b53-471b-956f-71701457972a.md) to the
# Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2          # xlUnderlineStyleSingle
$hyperlinkColor     = 15570276   # BGR for RGB(0x64,0x95,0xED) == theme hyperlink blue used by this workbook
$dateFormat         = "yyyy-mm-dd HH:mm:ss"

$fileName      = "01504b35-5b53-471b-956f-71701457972a.md"
$pathAndName   = "e2e\01504b35-5b53-471b-956f-71701457972a.md"
$extension     = ".md"
$status        = "Handed back: in sync with en-US"
$sourcePath    = "e2e"
$priority      = "ht"
$contentDup    = "True"
$zhXliff       = "01504b35-5b53-471b-956f-71701457972a.06f4bc6d2d107ccc8ecab03b060ed371ef11c6f8.zh-cn.xlf"
$zhHandoffDate  = "2016-09-03 08:48:15"
$zhHandbackDate = "2016-09-03 08:48:33"
$deXliff       = "01504b35-5b53-471b-956f-71701457972a.06f4bc6d2d107ccc8ecab03b060ed371ef11c6f8.de-de.xlf"
$deHandoffDate  = "2016-09-03 08:48:19"
$deHandbackDate = "2016-09-03 08:48:40"

# -------------------------------------------------------------------------
# Overview sheet
# -------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathAndName
$wsOverview.Range("C4").Value = $extension
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = $deHandoffDate
$wsOverview.Range("G4").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a3c4d5e6f708192a3b4c5d6e7f8091a2b3c4d5e/e2e/01504b35-5b53-471b-956f-71701457972a.md", "", "", $pathAndName)
$wsOverview.Range("B4").Font.Underline = $hyperlinkUnderline
$wsOverview.Range("B4").Font.Color = $hyperlinkColor

# -------------------------------------------------------------------------
# zh-cn sheet
# -------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = $extension
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = $sourcePath
$wsZh.Range("E4").Value = $priority
$wsZh.Range("F4").Value = $contentDup
$wsZh.Range("G4").Value = $zhXliff
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("H4").NumberFormat = $dateFormat
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXliff
$wsZh.Range("K4").Value = $zhHandbackDate
$wsZh.Range("K4").NumberFormat = $dateFormat
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a3c4d5e6f708192a3b4c5d6e7f8091a2b3c4d5e/e2e/01504b35-5b53-471b-956f-71701457972a.md", "", "", $fileName)
$wsZh.Range("A4").Font.Underline = $hyperlinkUnderline
$wsZh.Range("A4").Font.Color = $hyperlinkColor

$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3b4c5d6e7f8091a2b3c4d5e6f708192a3b4c5d6/e2e/01504b35-5b53-471b-956f-71701457972a.md", "", "", $fileName)
$wsZh.Range("I4").Font.Underline = $hyperlinkUnderline
$wsZh.Range("I4").Font.Color = $hyperlinkColor

# -------------------------------------------------------------------------
# de-de sheet
# -------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = $extension
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = $sourcePath
$wsDe.Range("E4").Value = $priority
$wsDe.Range("F4").Value = $contentDup
$wsDe.Range("G4").Value = $deXliff
$wsDe.Range("H4").Value = $deHandoffDate
$wsDe.Range("H4").NumberFormat = $dateFormat
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXliff
$wsDe.Range("K4").Value = $deHandbackDate
$wsDe.Range("K4").NumberFormat = $dateFormat
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a3c4d5e6f708192a3b4c5d6e7f8091a2b3c4d5e/e2e/01504b35-5b53-471b-956f-71701457972a.md", "", "", $fileName)
$wsDe.Range("A4").Font.Underline = $hyperlinkUnderline
$wsDe.Range("A4").Font.Color = $hyperlinkColor

$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4c5d6e7f8091a2b3c4d5e6f708192a3b4c5d6e7/e2e/01504b35-5b53-471b-956f-71701457972a.md", "", "", $fileName)
$wsDe.Range("I4").Font.Underline = $hyperlinkUnderline
$wsDe.Range("I4").Font.Color = $hyperlinkColor

Write-Host "Applied handback report update."
